$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.389.32'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.841.93'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6253'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07386'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2890'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.85'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07719'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '1.838.84'
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.964'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6712'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001026'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.275'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').Value = '29.390.56'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '234.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.283'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.99%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '157.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.475'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1343'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.07297'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.10%  '
$ws.Range('E29').Value = '  +5.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.474'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.039'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('E32').Value = '  -1.47%  '
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7127'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.581'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01833'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.74%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.233.73'
$ws.Range('E38').Value = '  -2.18%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.778'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.790'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9537'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.30%  '
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = '1.993.55'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.98'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000117'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.65%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.698'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.964'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.879'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3881'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.53%  '
